# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the zh-cn and
# de-de sheets' first data row (row 2) with new report-generation times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 02:43:00"
$wsZhCn.Range("H2").Value = "2016-03-14 02:43:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 02:43:03"
$wsDeDe.Range("H2").Value = "2016-03-14 02:43:27"
